$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the week label (shared string used by B9, and cascades via formulas
# to H9, B27, H27, B43). "SEMANA 17 ..." -> "SEMANA 18 ..."
$ws.Range("B9").Value = "SEMANA   18  DEL    02      Al   01   DE   MAYO          2022"

# Bonus/"Extra" amount for this period: 1400 -> 1680 (K24 SUM recalculates).
$ws.Range("K21").Value = 1680

# Second block "extra" line item: 1250 -> 0 (E41 SUM recalculates).
$ws.Range("E40").Value = 0

# Closing date used throughout the sheet (TODAY()) advances from
# 2022-04-30 to 2022-05-07 (cierre 7 may 22). The dependent cells I14,
# C32, I32 and C48 all reference C14 (directly or transitively) and pick
# up the new value automatically.
$ws.Range("C14").Value = 44688

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("D40").Select() | Out-Null
